$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the quarter reference date as plain text ("01/07/2024").
# Writing a date-looking string straight into .Value gets auto-parsed into
# a real date serial by Excel, so briefly force Text format, assign the
# literal string, then clear the format again (mirrors typing into a
# Text-formatted column then resetting formatting - value stays text,
# style reverts to the sheet default).
$dateRng = $ws.Range("C2:C10")
$dateRng.NumberFormat = "@"
$dateRng.Value = "01/10/2024"
$dateRng.ClearFormats()

# Row 2: Santa Catarina - updated value for new quarter
$ws.Range("D2").Value = 56.02

# Row 3: São Paulo - updated value for new quarter
$ws.Range("D3").Value = 55.75

# Row 4: was Goiás -> now Rio Grande do Sul (reordered), updated value
$ws.Range("A4").Value = "Rio Grande do Sul"
$ws.Range("D4").Value = 54.98

# Row 5: was Rio Grande do Sul -> now Distrito Federal (reordered), updated value
$ws.Range("A5").Value = "Distrito Federal"
$ws.Range("D5").Value = 54.88

# Row 6: Mato Grosso - updated value
$ws.Range("D6").Value = 54.5

# Row 7: was Distrito Federal -> now Goiás (reordered), updated value
$ws.Range("A7").Value = "Goiás"
$ws.Range("D7").Value = 54.02

# Row 8: Sergipe - updated value and ranking
$ws.Range("D8").Value = 46.4
$ws.Range("E8").Value = "18º"

# Row 9: Brasil - updated value
$ws.Range("D9").Value = 50.87

# Row 10: Nordeste - updated value
$ws.Range("D10").Value = 44.25
